# Insert a new weekly price record as row 730, shifting all subsequent
# rows (old 730-787) down by one (new 731-788).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(730).Insert()

$ws.Range("A730").Value = 6
$ws.Range("B730").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C730").Value = "Metropolitana"
$ws.Range("D730").Value = 45223
$ws.Range("E730").Value = 13
$ws.Range("F730").Value = 100112039
$ws.Range("G730").Value = "Ciboulette"
$ws.Range("H730").Value = "Sin especificar"
$ws.Range("I730").Value = "Primera"
$ws.Range("J730").Value = 580
$ws.Range("K730").Value = 1200
$ws.Range("L730").Value = 1300
$ws.Range("M730").Value = 1257
$ws.Range("N730").Value = "`$/docena de atados"
$ws.Range("O730").Value = "Región Metropolitana"
$ws.Range("P730").Value = 419
$ws.Range("Q730").Value = 3
$ws.Range("R730").Value = "Hortaliza"
